$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.048564009522522
$ws.Range("C2").Value = 0.3391783373210728
$ws.Range("E2").Value = 0.03282930602375345
$ws.Range("F2").Value = 0.4443680307746263
$ws.Range("G2").Value = 0.002541578262752807
$ws.Range("I2").Value = 2.271183132967721
$ws.Range("L2").Value = 0.2975070824688686
$ws.Range("B3").Value = 1.947589550913108
$ws.Range("C3").Value = 0.301067710024455
$ws.Range("E3").Value = 0.03288561977588311
$ws.Range("F3").Value = 0.387822817061874
$ws.Range("G3").Value = 0.002548195112812574
$ws.Range("I3").Value = 2.20423224319596
$ws.Range("L3").Value = 0.2871410443532199
$ws.Range("B4").Value = 1.887176263788206
$ws.Range("C4").Value = 0.2778217906077884
$ws.Range("E4").Value = 0.03292420771249194
$ws.Range("F4").Value = 0.3531389305168915
$ws.Range("G4").Value = 0.002552464330100363
$ws.Range("I4").Value = 2.163953270228461
$ws.Range("L4").Value = 0.2809706534582119
$ws.Range("B5").Value = 1.862953556297043
$ws.Range("C5").Value = 0.268386198267109
$ws.Range("E5").Value = 0.03294094160173988
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002554256185162608
$ws.Range("I5").Value = 2.147743181767297
$ws.Range("L5").Value = 0.2785047007317587
$ws.Range("B6").Value = 1.858955252930627
$ws.Range("C6").Value = 0.2668216365477178
$ws.Range("E6").Value = 0.03294378120529029
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002554556874953036
$ws.Range("I6").Value = 2.145063691776073
$ws.Range("L6").Value = 0.2780981522342358
$ws.Range("B7").Value = 1.88684798689269
$ws.Range("C7").Value = 0.2776943897587216
$ws.Range("E7").Value = 0.03292442930562289
$ws.Range("F7").Value = 0.3529483938368969
$ws.Range("G7").Value = 0.002552488284521221
$ws.Range("I7").Value = 2.163733835672915
$ws.Range("L7").Value = 0.2809372007120459
$ws.Range("B8").Value = 2.013417407121153
$ws.Range("C8").Value = 0.3260049036360897
$ws.Range("E8").Value = 0.03284789069824723
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("G8").Value = 0.002543817026647478
$ws.Range("I8").Value = 2.247923598704659
$ws.Range("L8").Value = 0.2938923076201689
$ws.Range("B9").Value = 2.274330584480992
$ws.Range("C9").Value = 0.4220358820753631
$ws.Range("E9").Value = 0.03272961675064545
$ws.Range("F9").Value = 0.5661985755042025
$ws.Range("G9").Value = 0.002528441524829806
$ws.Range("I9").Value = 2.419813178654266
$ws.Range("L9").Value = 0.3208574981389916
$ws.Range("B10").Value = 2.473982527900375
$ws.Range("C10").Value = 0.4934849064833884
$ws.Range("E10").Value = 0.03266212041972949
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002518125117484393
$ws.Range("I10").Value = 2.550551836813327
$ws.Range("L10").Value = 0.3416474848295934
$ws.Range("B11").Value = 2.566586249264446
$ws.Range("C11").Value = 0.5262070145349753
$ws.Range("E11").Value = 0.03263563058677632
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002513641918704845
$ws.Range("I11").Value = 2.61106278211696
$ws.Range("L11").Value = 0.3513241366049726
$ws.Range("B12").Value = 2.601912609044859
$ws.Range("C12").Value = 0.5386314294879071
$ws.Range("E12").Value = 0.03262620607654965
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002511974199931446
$ws.Range("I12").Value = 2.634131186256326
$ws.Range("L12").Value = 0.3550204225523714
$ws.Range("B13").Value = 2.594292858514905
$ws.Range("C13").Value = 0.5359541028497006
$ws.Range("E13").Value = 0.03262820882383921
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002512332042981846
$ws.Range("I13").Value = 2.629156052414601
$ws.Range("L13").Value = 0.35422293338641
$ws.Range("B14").Value = 2.569487353675356
$ws.Range("C14").Value = 0.5272285022142
$ws.Range("E14").Value = 0.03263484306554787
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002513504115007956
$ws.Range("I14").Value = 2.612957513577072
$ws.Range("L14").Value = 0.3516275894893113
$ws.Range("B15").Value = 2.554327137263954
$ws.Range("C15").Value = 0.5218882023024207
$ws.Range("E15").Value = 0.03263898574827961
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002514225940344121
$ws.Range("I15").Value = 2.603055676275432
$ws.Range("L15").Value = 0.3500420399209361
$ws.Range("B16").Value = 2.467966769430745
$ws.Range("C16").Value = 0.4913510018316174
$ws.Range("E16").Value = 0.03266393644583321
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002518422309884991
$ws.Range("I16").Value = 2.546618608545572
$ws.Range("L16").Value = 0.3410195398781184
$ws.Range("B17").Value = 2.41544591015753
$ws.Range("C17").Value = 0.4726748495634752
$ws.Range("E17").Value = 0.03268032270426957
$ws.Range("F17").Value = 0.6400460337215605
$ws.Range("G17").Value = 0.002521050238617537
$ws.Range("I17").Value = 2.512265553369929
$ws.Range("L17").Value = 0.3355409685672726
$ws.Range("B18").Value = 2.385404741833554
$ws.Range("C18").Value = 0.4619533865849803
$ws.Range("E18").Value = 0.03269014430142003
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002522581511461637
$ws.Range("I18").Value = 2.492603912059565
$ws.Range("L18").Value = 0.3324104433880279
$ws.Range("B19").Value = 2.375261980093967
$ws.Range("C19").Value = 0.4583267648181391
$ws.Range("E19").Value = 0.03269353783663242
$ws.Range("F19").Value = 0.6191636801734006
$ws.Range("G19").Value = 0.002523103373702958
$ws.Range("I19").Value = 2.485963380739861
$ws.Range("L19").Value = 0.3313540272181399
$ws.Range("B20").Value = 2.421019492798393
$ws.Range("C20").Value = 0.4746608172599167
$ws.Range("E20").Value = 0.03267853730323189
$ws.Range("F20").Value = 0.6429339538360921
$ws.Range("G20").Value = 0.002520768447671103
$ws.Range("I20").Value = 2.515912385202284
$ws.Range("L20").Value = 0.3361220363079696
$ws.Range("B21").Value = 2.576766267014534
$ws.Range("C21").Value = 0.5297905062855079
$ws.Range("E21").Value = 0.03263287795939851
$ws.Range("F21").Value = 0.7228739723492197
$ws.Range("G21").Value = 0.002513159037305976
$ws.Range("I21").Value = 2.617711190931004
$ws.Range("L21").Value = 0.3523890343666665
$ws.Range("B22").Value = 2.680069002963421
$ws.Range("C22").Value = 0.5660156328983703
$ws.Range("E22").Value = 0.03260657311338777
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.00250836045131551
$ws.Range("I22").Value = 2.685143908287699
$ws.Range("L22").Value = 0.3632068823265939
$ws.Range("B23").Value = 2.624794832215969
$ws.Range("C23").Value = 0.546663221131837
$ws.Range("E23").Value = 0.03262028872277667
$ws.Range("F23").Value = 0.7472568307916134
$ws.Range("G23").Value = 0.002510905635558547
$ws.Range("I23").Value = 2.649069641273371
$ws.Range("L23").Value = 0.3574159905568592
$ws.Range("B24").Value = 2.418499197359097
$ws.Range("C24").Value = 0.4737629124721821
$ws.Range("E24").Value = 0.03267934323436017
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002520895781511231
$ws.Range("I24").Value = 2.514263377605232
$ws.Range("L24").Value = 0.3358592758781924
$ws.Range("B25").Value = 2.202364544647025
$ws.Range("C25").Value = 0.3959075493743285
$ws.Range("E25").Value = 0.03275820686337871
$ws.Range("F25").Value = 0.5279251897347308
$ws.Range("G25").Value = 0.002532427973887575
$ws.Range("I25").Value = 2.372553312771942
$ws.Range("L25").Value = 0.3133928584521612
